$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert a brand-new column at N (shifts old N..AC one column right, to
#    O..AD) to make room for the new "localdb" command-type column.
# ---------------------------------------------------------------------------
$ws.Columns("N").Insert(-4161)

# ---------------------------------------------------------------------------
# 2. Populate the new column N with the "localdb" command type: a header in
#    row 1 and its six commands in rows 2-7.
# ---------------------------------------------------------------------------
$ws.Range("N1").Value2 = "localdb"
$ws.Range("N2").Value2 = "cloneTable(var,source,target)"
$ws.Range("N3").Value2 = "dropTables(var,tables)"
$ws.Range("N4").Value2 = "exportCSV(sql,output)"
$ws.Range("N5").Value2 = "importRecords(var,sourceDb,sql,table)"
$ws.Range("N6").Value2 = "purge(var)"
$ws.Range("N7").Value2 = "runSQLs(var,sqls)"

# ---------------------------------------------------------------------------
# 3. Insert a new row-14 entry in column A (the alphabetically sorted list of
#    "target" names) shifting A14:A29 down to A15:A30, then set the new
#    A14 cell to "localdb" (alphabetically between "json" and "macro").
# ---------------------------------------------------------------------------
$ws.Range("A14").Insert(-4121)
$ws.Range("A14").Value2 = "localdb"

# ---------------------------------------------------------------------------
# 4. Update every defined name whose target range shifted right by one
#    column because of the column insert in step 1.
# ---------------------------------------------------------------------------
$wb.Names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$wb.Names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$wb.Names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$wb.Names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$wb.Names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$wb.Names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$wb.Names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$wb.Names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$wb.Names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$wb.Names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$wb.Names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$wb.Names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"

# "target" grew by one row because of the row-14 insert in step 3.
$wb.Names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"

# ---------------------------------------------------------------------------
# 5. Register the new "localdb" defined name.
# ---------------------------------------------------------------------------
$wb.Names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
